# Insert a new data row at row 252, shifting existing rows 252-361 down to 253-362.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with its values.
$ws.Cells.Item(252, 1).Value = 7
$ws.Cells.Item(252, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(252, 3).Value = "Ñuble"
$ws.Cells.Item(252, 4).Value = 45141
$ws.Cells.Item(252, 5).Value = 16
$ws.Cells.Item(252, 6).Value = 100112043
$ws.Cells.Item(252, 7).Value = "Pepino ensalada"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 60
$ws.Cells.Item(252, 11).Value = 13000
$ws.Cells.Item(252, 12).Value = 13000
$ws.Cells.Item(252, 13).Value = 13000
$ws.Cells.Item(252, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(252, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(252, 16).Value = 217
$ws.Cells.Item(252, 17).Value = 60
$ws.Cells.Item(252, 18).Value = "Hortaliza"
